$wb = $excel.ActiveWorkbook

# --- Rename sheets (task-order run ids bumped) ---
$wb.Worksheets.Item(1).Name = "GNG_TO-16509962046071513"
$wb.Worksheets.Item(2).Name = "NB_TO-1650996206111124"
$wb.Worksheets.Item(3).Name = "RS_TO-1650996206111124"
$wb.Worksheets.Item(4).Name = "TOL_TO-16509962061591296"
$wb.Worksheets.Item(5).Name = "vSAT_TO-16509962062231512"

# --- Sheet 1: GNG ---
$ws = $wb.Worksheets.Item(1)
$ws.Range("B2").Value = "go_stims-16509962045751183.csv"
$ws.Range("B3").Value = "GNG_stims-1650996204591153.csv"
$ws.Range("B4").Value = "go_stims-1650996204591153.csv"
$ws.Range("B5").Value = "GNG_stims-16509962046071513.csv"

# --- Sheet 2: NB ---
$ws = $wb.Worksheets.Item(2)
$ws.Range("B2").Value = "ZB-match_4-16509962047991204.csv"
$ws.Range("B3").Value = "OB-16509962054951217.csv"
$ws.Range("B4").Value = "TB-16509962058231525.csv"
$ws.Range("B5").Value = "TB-16509962060951233.csv"
$ws.Range("B6").Value = "ZB-match_5-16509962046791196.csv"
$ws.Range("B7").Value = "OB-16509962051591196.csv"
$ws.Range("B8").Value = "TB-16509962060471582.csv"
$ws.Range("B9").Value = "ZB-match_2-16509962048311222.csv"
$ws.Range("B10").Value = "OB-16509962053111174.csv"

# --- Sheet 3: RS --- (no cell content changes)

# --- Sheet 4: TOL ---
$ws = $wb.Worksheets.Item(4)
$ws.Range("B2").Value = "MM_stims-16509962061271274.csv"
$ws.Range("B3").Value = "ZM_stims-1650996206111124.csv"
$ws.Range("B4").Value = "MM_stims-16509962061431558.csv"
$ws.Range("B5").Value = "ZM_stims-16509962061271274.csv"
$ws.Range("B6").Value = "MM_stims-16509962061591296.csv"
$ws.Range("B7").Value = "ZM_stims-16509962061431558.csv"

# --- Sheet 5: vSAT ---
$ws = $wb.Worksheets.Item(5)
$ws.Range("B2").Value = "SAT_stims-16509962061591296.csv"
$ws.Range("B3").Value = "vSAT_stims-16509962061911545.csv"
$ws.Range("B4").Value = "vSAT_stims-16509962062071204.csv"
$ws.Range("B5").Value = "SAT_stims-16509962061751244.csv"
